# Automatic update of files.
#
# The species-observation rows 10-15 and 17-23 were refreshed with a newer
# export of the same records (row 16 is unchanged). Most rows simply swap
# their full content with another row in the block (same underlying
# sighting, re-ordered), a couple of coordinate/time-only tweaks happen on
# rows whose species data stayed put, and row 20 ends up holding a brand
# new "Rosenticka" observation. Row 13 also picks up an "Aktivitet" (M)
# comment that used to live on row 21, which loses it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("A10").Value = 111577919
$ws.Range("B10").Value = 96348
$ws.Range("D10").Value = 'VU'
$ws.Range("E10").Value = 220787
$ws.Range("F10").Value = 'Knärot'
$ws.Range("G10").Value = 'Goodyera repens'
$ws.Range("H10").Value = '(L.) R. Br.'
$ws.Range("Q10").Value = 562836.2055113926
$ws.Range("R10").Value = 6954423.824987715
$ws.Range("Z10").Value = '00:00'
$ws.Range("AB10").Value = '00:00'

# Row 11
$ws.Range("A11").Value = 111577080
$ws.Range("B11").Value = 89405
$ws.Range("D11").Value = 'NT'
$ws.Range("E11").Value = 1202
$ws.Range("F11").Value = 'Ullticka'
$ws.Range("G11").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H11").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q11").Value = 562937.8525077751
$ws.Range("R11").Value = 6954467.524316943

# Row 12
$ws.Range("A12").Value = 111577273
$ws.Range("B12").Value = 89686
$ws.Range("E12").Value = 658
$ws.Range("F12").Value = 'Rosenticka'
$ws.Range("G12").Value = 'Rhodofomes roseus'
$ws.Range("H12").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q12").Value = 562824.7977144517
$ws.Range("R12").Value = 6954323.105396069
$ws.Range("Z12").Value = '17:24'
$ws.Range("AB12").Value = '17:24'

# Row 13
$ws.Range("A13").Value = 111577591
$ws.Range("B13").Value = 56543
$ws.Range("D13").Value = 'NT'
$ws.Range("E13").Value = 103021
$ws.Range("F13").Value = 'Talltita'
$ws.Range("G13").Value = 'Poecile montanus'
$ws.Range("H13").Value = '(Conrad von Baldenstein, 1827)'
$ws.Range("Q13").Value = 562822.1033927511
$ws.Range("R13").Value = 6954368.028004575

# Row 14
$ws.Range("A14").Value = 111578090
$ws.Range("B14").Value = 89405
$ws.Range("D14").Value = 'NT'
$ws.Range("E14").Value = 1202
$ws.Range("F14").Value = 'Ullticka'
$ws.Range("G14").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H14").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q14").Value = 562890.6102569005
$ws.Range("R14").Value = 6954486.814324431
$ws.Range("Z14").Value = '18:30'
$ws.Range("AB14").Value = '18:30'

# Row 15
$ws.Range("A15").Value = 111575408
$ws.Range("Q15").Value = 562539.034657649
$ws.Range("R15").Value = 6954609.073577877
$ws.Range("Z15").Value = '16:21'
$ws.Range("AB15").Value = '16:21'

# Row 16 is unchanged

# Row 17
$ws.Range("A17").Value = 111577743
$ws.Range("B17").Value = 96348
$ws.Range("D17").Value = 'VU'
$ws.Range("E17").Value = 220787
$ws.Range("F17").Value = 'Knärot'
$ws.Range("G17").Value = 'Goodyera repens'
$ws.Range("H17").Value = '(L.) R. Br.'
$ws.Range("Q17").Value = 562802.8660743404
$ws.Range("R17").Value = 6954388.771485241
$ws.Range("Z17").Value = '17:59'
$ws.Range("AB17").Value = '17:59'

# Row 18
$ws.Range("A18").Value = 111577611
$ws.Range("Q18").Value = 562810.9079359611
$ws.Range("R18").Value = 6954400.856378952

# Row 19
$ws.Range("A19").Value = 111577347
$ws.Range("Q19").Value = 562796.503171768
$ws.Range("R19").Value = 6954336.792844097
$ws.Range("Z19").Value = '00:00'
$ws.Range("AB19").Value = '00:00'

# Row 20
$ws.Range("A20").Value = 111577994
$ws.Range("B20").Value = 89686
$ws.Range("D20").Value = 'NT'
$ws.Range("E20").Value = 658
$ws.Range("F20").Value = 'Rosenticka'
$ws.Range("G20").Value = 'Rhodofomes roseus'
$ws.Range("H20").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q20").Value = 562890.6102569005
$ws.Range("R20").Value = 6954486.814324431

# Row 21
$ws.Range("A21").Value = 111577838
$ws.Range("B21").Value = 96348
$ws.Range("D21").Value = 'VU'
$ws.Range("E21").Value = 220787
$ws.Range("F21").Value = 'Knärot'
$ws.Range("G21").Value = 'Goodyera repens'
$ws.Range("H21").Value = '(L.) R. Br.'
$ws.Range("Q21").Value = 562801.9298062191
$ws.Range("R21").Value = 6954389.67147268
$ws.Range("Z21").Value = '18:14'
$ws.Range("AB21").Value = '18:14'

# Row 22
$ws.Range("A22").Value = 111577804
$ws.Range("Q22").Value = 562814.775380839
$ws.Range("R22").Value = 6954390.834027934
$ws.Range("Z22").Value = '00:00'
$ws.Range("AB22").Value = '00:00'

# Row 23
$ws.Range("A23").Value = 111577193
$ws.Range("B23").Value = 89845
$ws.Range("D23").Value = 'VU'
$ws.Range("E23").Value = 1209
$ws.Range("F23").Value = 'Rynkskinn'
$ws.Range("G23").Value = 'Phlebia centrifuga'
$ws.Range("H23").Value = 'P.Karst.'
$ws.Range("Q23").Value = 562885.3077477051
$ws.Range("R23").Value = 6954427.514711756

# Row 13 gains an M13 cell (Aktivitet) that was previously on row 21
$ws.Range("M13").Value = 'lockläte, övriga läten'

# Row 21 loses its M21 cell (Aktivitet) -- cleared so the cell becomes empty
$ws.Range("M21").Value = ""

